$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.122.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("E2").Style = "Normal"

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.251.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E3").Style = "Normal"

# Row 4 - TetherUSD
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("E4").Style = "Normal"

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("E5").Style = "Normal"

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("E6").Style = "Normal"

# Row 7 - USDC
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E7").Style = "Normal"

# Row 8 - LidoStakedEther
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.248.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("E8").Style = "Normal"

# Row 9 - XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("E9").Style = "Normal"

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("E10").Style = "Normal"

# Row 11 - Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("E11").Style = "Normal"

# Row 12 - Cardano
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.13%  "
$ws.Range("E12").Style = "Normal"

# Row 13 - ShibaInu
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000267"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("E13").Style = "Normal"

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("E14").Style = "Normal"

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.784.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("E15").Style = "Normal"

# Row 16 - WrappedBTC
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.171.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("E16").Style = "Normal"

# Row 17 - WrappedEther
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.250.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("E17").Style = "Normal"

# Row 18 - Polkadot
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.90%  "
$ws.Range("E18").Style = "Normal"

# Row 19 - TRON
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("E19").Style = "Normal"

# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "496.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("E20").Style = "Normal"

# Row 21 - Chainlink
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("E21").Style = "Normal"

# Row 22 - Polygon
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.742"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("E22").Style = "Normal"

# Row 23 - Uniswap
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.44%  "
$ws.Range("E23").Style = "Normal"

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.71%  "
$ws.Range("E24").Style = "Normal"

# Row 25 - Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.27%  "
$ws.Range("E25").Style = "Normal"

# Row 26 - Dai
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E26").Style = "Normal"

# Row 27 - PancakeSwap
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E27").Style = "Normal"

# Row 28 - RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.38%  "
$ws.Range("E28").Style = "Normal"

# Row 29 - ImmutableX->Hedera
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Hedera"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.132"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +43.47%  "
$ws.Range("E29").Style = "Normal"

# Row 30 - Hedera->ImmutableX
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("E30").Style = "Normal"

# Row 31 - NEARProtocol
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("E31").Style = "Normal"

# Row 32 - Stacks
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -7.58%  "
$ws.Range("E32").Style = "Normal"

# Row 33 - EthereumClassic
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.99%  "
$ws.Range("E33").Style = "Normal"

# Row 34 - FirstDigitalUSD
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E34").Style = "Normal"

# Row 35 - Mantle
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.08%  "
$ws.Range("E35").Style = "Normal"

# Row 36 - Filecoin
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.75%  "
$ws.Range("E36").Style = "Normal"

# Row 37 - dogwifhat
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +12.89%  "
$ws.Range("E37").Style = "Normal"

# Row 38 - OKB
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.69"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("E38").Style = "Normal"

# Row 39 - Bittensor
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "490.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.15%  "
$ws.Range("E39").Style = "Normal"

# Row 40 - PEPE
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0773"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("E40").Style = "Normal"

# Row 41 - VeChain
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0419"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E41").Style = "Normal"

# Row 42 - Kaspa
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("E42").Style = "Normal"

# Row 43 - Cosmos
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("E43").Style = "Normal"

# Row 44 - Fetch.AI
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E44").Style = "Normal"

# Row 45 - Maker
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.992.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.75%  "
$ws.Range("E45").Style = "Normal"

# Row 46 - TheGraph
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.70%  "
$ws.Range("E46").Style = "Normal"

# Row 47 - InjectiveProtocol
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("E47").Style = "Normal"

# Row 48 - ThetaToken
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("E48").Style = "Normal"

# Row 49 - Stellar
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("E49").Style = "Normal"

# Row 51 - Monero
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.38%  "
$ws.Range("E51").Style = "Normal"
